$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text even when $value looks numeric,
    # mirroring the original inline-string cells (no lasting style change).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue $ws.Range("D2") "41.076.59"
Set-TextValue $ws.Range("D3") "2.174.94"
Set-TextValue $ws.Range("D5") "251.66"
Set-TextValue $ws.Range("D7") "66.03"
Set-TextValue $ws.Range("D9") "0.575"
Set-TextValue $ws.Range("D10") "58.97"
Set-TextValue $ws.Range("D11") "36.40"
Set-TextValue $ws.Range("D12") "0.0933"
Set-TextValue $ws.Range("D13") "0.104"
Set-TextValue $ws.Range("D14") "6.83"
Set-TextValue $ws.Range("D15") "2.498.49"
Set-TextValue $ws.Range("D16") "14.24"
Set-TextValue $ws.Range("D17") "0.844"
Set-TextValue $ws.Range("D18") "2.167.02"
Set-TextValue $ws.Range("D19") "41.057.33"
Set-TextValue $ws.Range("D20") "0.0₃0945"
Set-TextValue $ws.Range("D21") "71.50"
Set-TextValue $ws.Range("D22") "6.04"
Set-TextValue $ws.Range("D23") "229.95"
Set-TextValue $ws.Range("D27") "11.33"
Set-TextValue $ws.Range("D28") "2.40"
Set-TextValue $ws.Range("D30") "167.89"
Set-TextValue $ws.Range("D31") "20.14"
Set-TextValue $ws.Range("D34") "0.0747"
Set-TextValue $ws.Range("D36") "4.50"
Set-TextValue $ws.Range("D37") "3.92"
Set-TextValue $ws.Range("D38") "24.45"
Set-TextValue $ws.Range("D41") "2.21"
Set-TextValue $ws.Range("D43") "60.89"
Set-TextValue $ws.Range("D44") "11.31"
Set-TextValue $ws.Range("D45") "8.46"
Set-TextValue $ws.Range("D48") "0.187"
Set-TextValue $ws.Range("D49") "1.14"

# --- Volume(1h) (column E) updates ---
Set-TextValue $ws.Range("E3") "  -2.33%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("E5") "  +0.57%  "
Set-TextValue $ws.Range("E6") "  -3.23%  "
Set-TextValue $ws.Range("E7") "  -8.10%  "
Set-TextValue $ws.Range("E8") "  +0.08%  "
Set-TextValue $ws.Range("E9") "  -4.71%  "
Set-TextValue $ws.Range("E10") "  +0.87%  "
Set-TextValue $ws.Range("E11") "  -11.39%  "
Set-TextValue $ws.Range("E12") "  -3.78%  "
Set-TextValue $ws.Range("E13") "  -1.04%  "
Set-TextValue $ws.Range("E14") "  -5.55%  "
Set-TextValue $ws.Range("E15") "  -2.31%  "
Set-TextValue $ws.Range("E16") "  -5.03%  "
Set-TextValue $ws.Range("E17") "  -2.78%  "
Set-TextValue $ws.Range("E18") "  -2.33%  "
Set-TextValue $ws.Range("E19") "  -1.60%  "
Set-TextValue $ws.Range("E20") "  -1.91%  "
Set-TextValue $ws.Range("E21") "  -2.01%  "
Set-TextValue $ws.Range("E23") "  -2.18%  "
Set-TextValue $ws.Range("E24") "  -3.97%  "
Set-TextValue $ws.Range("E25") "  -4.39%  "
Set-TextValue $ws.Range("E27") "  +5.18%  "
Set-TextValue $ws.Range("E28") "  -5.21%  "
Set-TextValue $ws.Range("E29") "  +0.08%  "
Set-TextValue $ws.Range("E30") "  -1.93%  "
Set-TextValue $ws.Range("E31") "  -3.25%  "
Set-TextValue $ws.Range("E32") "  -3.12%  "
Set-TextValue $ws.Range("E33") "  +1.00%  "
Set-TextValue $ws.Range("E34") "  +1.92%  "
Set-TextValue $ws.Range("E35") "  -2.81%  "
Set-TextValue $ws.Range("E36") "  -4.99%  "
Set-TextValue $ws.Range("E37") "  -2.31%  "
Set-TextValue $ws.Range("E38") "  -5.25%  "
Set-TextValue $ws.Range("E39") "  +0.41%  "
Set-TextValue $ws.Range("E40") "  +12.94%  "
Set-TextValue $ws.Range("E41") "  -3.89%  "
Set-TextValue $ws.Range("E42") "  -8.29%  "
Set-TextValue $ws.Range("E43") "  -8.51%  "
Set-TextValue $ws.Range("E44") "  -7.43%  "
Set-TextValue $ws.Range("E45") "  -3.58%  "
Set-TextValue $ws.Range("E48") "  -7.87%  "
Set-TextValue $ws.Range("E49") "  -3.09%  "

# --- Rows 46/47 swap: BinanceUSD/Cronos order flips with new values ---
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.0990"
Set-TextValue $ws.Range("E46") "  -3.28%  "

$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D47") "1.00"
Set-TextValue $ws.Range("E47") "  -0.23%  "

# --- Rows 50/51 swap: SynthetixNetwork/TrustWalletToken order flips with new values ---
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D50") "1.14"
Set-TextValue $ws.Range("E50") "  -4.23%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D51") "4.20"
Set-TextValue $ws.Range("E51") "  -10.05%  "
